# Estonia Meistriliiga workbook update (06-04-2024 01:36)
#
# This script reproduces a data refresh that:
#   1) Re-ordered several pairs of fixtures that share the exact same
#      kick-off date/time (the upstream feed re-sorted same-timestamp
#      matches), swapping every column except the running-rank column A.
#   2) Appended two brand-new fixture rows (129, 130) at the bottom of
#      the sheet - one finished match and one match that hasn't kicked
#      off yet (so it has no FTHG/FTAG/FTR or PL_AhOver/PL_AhUnder yet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 29   # column AC

function Swap-RowContent($row1, $row2) {
    # Swap every column from B through AC between two rows, leaving
    # column A (the running rank) untouched in both rows.
    for ($col = 2; $col -le $lastCol; $col++) {
        $cell1 = $ws.Cells.Item($row1, $col)
        $cell2 = $ws.Cells.Item($row2, $col)
        $val1 = $cell1.Value()
        $val2 = $cell2.Value()
        $cell1.Value = $val2
        $cell2.Value = $val1
    }
}

# Same-kick-off-time fixture pairs whose rows got swapped in this refresh.
Swap-RowContent 4 5
Swap-RowContent 64 65
Swap-RowContent 104 107
Swap-RowContent 105 106
Swap-RowContent 115 116
Swap-RowContent 120 121

# --- New fixture rows appended at the bottom of the sheet -----------------

# Row 129: JK Nomme United 2-1 JK Tammeka Tartu (played)
$ws.Cells.Item(129, 1).Value = 127
$ws.Cells.Item(129, 2).Value = 7721011
$ws.Cells.Item(129, 3).Value = "Estonia Meistriliiga"
$ws.Cells.Item(129, 4).Value = "Estonia Meistriliiga"
$ws.Cells.Item(129, 5).Value = 45387.54166666666
$ws.Cells.Item(129, 6).Value = "JK Nomme United"
$ws.Cells.Item(129, 7).Value = "JK Tammeka Tartu"
$ws.Cells.Item(129, 8).Value = 2
$ws.Cells.Item(129, 9).Value = 1
$ws.Cells.Item(129, 10).Value = "H"
$ws.Cells.Item(129, 11).Value = 3.8
$ws.Cells.Item(129, 12).Value = 3.5
$ws.Cells.Item(129, 13).Value = 1.85
$ws.Cells.Item(129, 14).Value = 4.2
$ws.Cells.Item(129, 15).Value = 3.4
$ws.Cells.Item(129, 16).Value = 1.8
$ws.Cells.Item(129, 17).Value = 0.75
$ws.Cells.Item(129, 18).Value = 1.8
$ws.Cells.Item(129, 19).Value = 2
$ws.Cells.Item(129, 20).Value = 2.5
$ws.Cells.Item(129, 21).Value = 1.9
$ws.Cells.Item(129, 22).Value = 1.9
$ws.Cells.Item(129, 23).Value = 3.2
$ws.Cells.Item(129, 24).Value = -1
$ws.Cells.Item(129, 25).Value = -1
$ws.Cells.Item(129, 26).Value = 0.8
$ws.Cells.Item(129, 27).Value = -1
$ws.Cells.Item(129, 28).Value = 0.8999999999999999
$ws.Cells.Item(129, 29).Value = -1

# Row 130: JK Nomme Kalju vs FC Levadia Tallinn (not played yet - no
# FTHG/FTAG/FTR and no closing PL_AhOver/PL_AhUnder).
$ws.Cells.Item(130, 1).Value = 128
$ws.Cells.Item(130, 2).Value = 7719647
$ws.Cells.Item(130, 3).Value = "Estonia Meistriliiga"
$ws.Cells.Item(130, 4).Value = "Estonia Meistriliiga"
$ws.Cells.Item(130, 5).Value = 45388.27083333334
$ws.Cells.Item(130, 6).Value = "JK Nomme Kalju"
$ws.Cells.Item(130, 7).Value = "FC Levadia Tallinn"
$ws.Cells.Item(130, 11).Value = 3.9
$ws.Cells.Item(130, 12).Value = 3.4
$ws.Cells.Item(130, 13).Value = 1.85
$ws.Cells.Item(130, 14).Value = 7.5
$ws.Cells.Item(130, 15).Value = 3.75
$ws.Cells.Item(130, 16).Value = 1.4
$ws.Cells.Item(130, 17).Value = 1.25
$ws.Cells.Item(130, 18).Value = 1.9
$ws.Cells.Item(130, 19).Value = 1.9
$ws.Cells.Item(130, 20).Value = 2.5
$ws.Cells.Item(130, 21).Value = 1.85
$ws.Cells.Item(130, 22).Value = 1.95
$ws.Cells.Item(130, 23).Value = 0
$ws.Cells.Item(130, 24).Value = 0
$ws.Cells.Item(130, 25).Value = 0
$ws.Cells.Item(130, 26).Value = 0
$ws.Cells.Item(130, 27).Value = 0

# Match the existing formatting used for every other data row: column A is
# bold/bordered/centered, column E carries the custom date/time number
# format. Copy it from row 128 (the previous last row) rather than
# hand-building style indices.
$ws.Range("A128").Copy()
$ws.Range("A129").PasteSpecial(-4122)
$ws.Range("A128").Copy()
$ws.Range("A130").PasteSpecial(-4122)
$ws.Range("E128").Copy()
$ws.Range("E129").PasteSpecial(-4122)
$ws.Range("E128").Copy()
$ws.Range("E130").PasteSpecial(-4122)
